# Fill in the quantity ("খাতা/পত্রের সংখ্যা") column G for the line items
# that previously had no entry. This drives the per-row rate*qty formulas
# in column I, and the grand-total SUM in I32 recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
$ws.Range("G29").Value = 1
